$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.741.57"
$ws.Range("E2").Value = "  +6.18%  "

$ws.Range("D3").Value = "3.054.49"
$ws.Range("E3").Value = "  +5.58%  "

$ws.Range("D4").Formula = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Formula = "'556.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.60%  "

$ws.Range("D6").Formula = "'142.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.46%  "

$ws.Range("D7").Formula = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.050.41"
$ws.Range("E8").Value = "  +5.48%  "

$ws.Range("D9").Formula = "'0.506"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.38%  "

$ws.Range("E10").Value = "  +10.43%  "

$ws.Range("D11").Formula = "'6.09"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.81%  "

$ws.Range("D12").Formula = "'0.477"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +11.84%  "

$ws.Range("E13").Value = "  +9.79%  "

$ws.Range("D14").Formula = "'34.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.48%  "

$ws.Range("D15").Value = "3.552.65"
$ws.Range("E15").Value = "  +4.67%  "

$ws.Range("D16").Value = "63.791.63"
$ws.Range("E16").Value = "  +6.10%  "

$ws.Range("E17").Value = "  +4.09%  "

$ws.Range("D18").Value = "3.052.73"
$ws.Range("E18").Value = "  +5.25%  "

$ws.Range("E19").Value = "  +6.37%  "

$ws.Range("D20").Formula = "'478.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.43%  "

$ws.Range("E21").Value = "  +9.04%  "

$ws.Range("E22").Value = "  +8.39%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Formula = "'14.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +20.55%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Formula = "'7.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +11.43%  "

$ws.Range("D25").Formula = "'81.35"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.12%  "

$ws.Range("D26").Formula = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +7.09%  "

$ws.Range("D28").Formula = "'7.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.63%  "

$ws.Range("D29").Formula = "'2.04"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("D30").Formula = "'0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Formula = "'26.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.88%  "

$ws.Range("E32").Value = "  +4.61%  "

$ws.Range("E33").Value = "  +9.93%  "

$ws.Range("D34").Formula = "'5.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.90%  "

$ws.Range("D35").Formula = "'6.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.55%  "

$ws.Range("D36").Formula = "'54.92"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.67%  "

$ws.Range("D37").Formula = "'0.0410"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.97%  "

$ws.Range("D38").Formula = "'446.68"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.69%  "

$ws.Range("E39").Value = "  +5.31%  "

$ws.Range("E40").Value = "  +24.48%  "

$ws.Range("D41").Value = "2.959.54"
$ws.Range("E41").Value = "  +3.97%  "

$ws.Range("E42").Value = "  +7.30%  "

$ws.Range("D43").Formula = "'0.113"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").Formula = "'27.81"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.16%  "

$ws.Range("D45").Formula = "'0.261"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +10.46%  "

$ws.Range("E46").Value = "  +14.77%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Formula = "'0.113"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.80%  "

$ws.Range("D49").Value = "0.0₃0515"
$ws.Range("E49").Value = "  +10.94%  "

$ws.Range("D50").Formula = "'116.85"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.99%  "

$ws.Range("D51").Formula = "'2.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.59%  "
